$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet1 (Sedan_HambaLG_f) hardpoint values ---
$ws1.Range("F5").Value2 = -0.0026557142857142869
$ws1.Range("F5").NumberFormat = "0.000"
$ws1.Range("G5").Value2 = 0.62
$ws1.Range("G5").NumberFormat = "0.00"
$ws1.Range("H5").Value2 = 0.65
$ws1.Range("H5").NumberFormat = "0.00"

$ws1.Range("F6").Value2 = 0.055166428571428582
$ws1.Range("F6").NumberFormat = "0.000"
$ws1.Range("G6").Value2 = 0.85
$ws1.Range("G6").NumberFormat = "0.00"
$ws1.Range("H6").Value2 = 0.19
$ws1.Range("H6").NumberFormat = "0.00"

# --- Sheet2 (Sedan_HambaLG_r) hardpoint values ---
$ws2.Range("F5").Value2 = 0.0026557142857142869
$ws2.Range("F5").NumberFormat = "0.000"
$ws2.Range("G5").Value2 = 0.62
$ws2.Range("G5").NumberFormat = "0.00"
$ws2.Range("H5").Value2 = 0.65
$ws2.Range("H5").NumberFormat = "0.00"

$ws2.Range("F6").Value2 = -0.055166428571428582
$ws2.Range("F6").NumberFormat = "0.000"
$ws2.Range("G6").Value2 = 0.85
$ws2.Range("G6").NumberFormat = "0.00"
$ws2.Range("H6").Value2 = 0.19
$ws2.Range("H6").NumberFormat = "0.00"

# --- Tab colors (theme Accent5, tint -0.25 ~ dark blue) ---
$ws1.Tab.Color = 10515524
$ws2.Tab.Color = 10515524

# --- View state: sheet1 becomes active/selected tab, frozen pane scrolled to column D ---
$ws2.Range("F25").Select()
$ws1.Activate()
$excel.ActiveWindow.ScrollColumn = 4
$ws1.Range("F34").Select()

Write-Output "done"
